$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'cold weather workout clothes for men',
    'cold weather workout pants for men',
    'cold weather workout women',
    'cold workout gear',
    'coloful koala leggings',
    'color run apparel',
    'colorfluka leggings',
    'colorful compression leggings women',
    'colorful koala leggings',
    'colorfulkoala legging',
    'colorfulkoala leggings with pockets',
    'coloruflkoala leggings',
    'comfortable leggings for women',
    'comfortable leggings women',
    'comfortable tights for women',
    'comfy pot',
    'coming over exo',
    'commando tights',
    'comment below',
    'common culture apparel',
    'common fig',
    'compreasion leggings',
    'compresion leggings',
    'compresion leggings for women',
    'compresion legs',
    'compresion short',
    'compresion shorts women',
    'compresion tights',
    'compress knee brace',
    'compress pants women',
    'compress running shorts women',
    'compressed leggings for women',
    'compressed shorts for women',
    'compressed tights',
    'compressin tank',
    'compressing vest',
    'compressio garment',
    'compression 2xu',
    'compression 2xu men',
    'compression 2xu women',
    'compression activewear for women',
    'compression adjustable breathable basketball protection',
    'compression apparel for women',
    'compression athletic leggings women',
    'compression back support',
    'compression back support women',
    'compression bands for legs',
    'compression basketball leggings',
    'compression basketball pants',
    'compression below knee',
    'compression bike shorts',
    'compression bike shorts womens high waist',
    'compression biker shorts for women high waist',
    'compression body women',
    'compression calf 2xu',
    'compression calf leggings',
    'compression calf running',
    'compression calves',
    'compression calves women',
    'compression capri',
    'compression capri leggings',
    'compression capri leggings women',
    'compression capri men',
    'compression capri men pack',
    'compression capri pants men',
    'compression capri women',
    'compression capris',
    'compression capris for men',
    'compression capris for women',
    'compression capris men',
    'compression capris women',
    'compression capris women high waist',
    'compression capris women with pocket',
    'compression clothes',
    'compression clothes for boys',
    'compression clothes women',
    'compression clothing',
    'compression clothing women',
    'compression cold gear',
    'compression cold weather leggings',
    'compression cord set',
    'compression crop top',
    'compression device',
    'compression devices legs',
    'compression dress pants women',
    'compression dry fit tank tops women',
    'compression exercise pants women',
    'compression fabric',
    'compression fitness pants',
    'compression for 20-30 varicose stockings',
    'compression for abdomen women',
    'compression for calf',
    'compression for calves',
    'compression for hamstring',
    'compression for knee',
    'compression for knees',
    'compression for knees women',
    'compression for legs',
    'compression for running',
    'compression for shins'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
